$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 16: add D16 = "Y"
$ws.Range("D16").Value = "Y"

# Row 22: add D22 = "TODO"
$ws.Range("D22").Value = "TODO"
$ws.Rows.Item(22).AutoFit()

# Update the view: scroll so A4 is the top-left visible cell, and select E16
$ws.Range("E16").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
